$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('prevention_vaccine', 'prevention_medicine', 35),
    @('prevention_medicine', 'prevention_vaccine', 20),
    @('vocative_thank_you', 'comment_positive', 10),
    @('comment_positive', 'vocative_thank_you', 9),
    @('test_who', 'test_virus', 9),
    @('vocative_help', 'bot_capabilities', 6),
    @('comment_offense', 'comment_negative', 4),
    @('prevention_home', 'prevention_medical_attention', 4),
    @('user_no_further_questions', 'vocative_thank_you', 4),
    @('covid_origins', 'covid_duration', 3),
    @('cc_philosophical', 'cc_religion', 3),
    @('bot_residence', 'bot_origin', 2),
    @('comment_negative', 'vocative_sorry', 2),
    @('covid_info', 'covid_duration', 2),
    @('covid_current_statistics', 'covid_info', 2),
    @('quarantine_general', 'quaratine_how_it_works', 2),
    @('covid_treatment', 'prevention_medicine', 2),
    @('bot_sexual', 'bot_personal_questions', 2),
    @('vocative_call', 'cc_lets_talk', 2),
    @('cc_religion', 'bot_personal_questions', 2),
    @('prevention_medical_attention', 'covid_symptoms', 2),
    @('bot_residence', 'bot_personal_questions', 2),
    @('cc_philosophical', 'bot_goal', 2),
    @('covid_info', 'covid_origins', 2),
    @('comment_smart', 'bot_personality', 2),
    @('coronavirus_info', 'covid_situation_last_update', 1),
    @('mask_general', 'mask_use_put', 1),
    @('covid_dangerous', 'spread_risk', 1),
    @('user_love', 'cc_philosophical', 1),
    @('bot_sing', 'bot_personal_questions', 1),
    @('user_particles', 'comment_negative', 1),
    @('comment_offense', 'comment_racist', 1),
    @('bot_appearance', 'bot_sexual', 1),
    @('prevention_home', 'comment_negative', 1),
    @('vocative_you_welcome', 'vocative_no', 1),
    @('covid_symptoms', 'prevention_medical_attention', 1),
    @('greeting_hello', 'start', 1),
    @('greeting_hello', 'greeting_goodbye', 1),
    @('bot_appearance', 'bot_personal_questions', 1),
    @('bot_music', 'bot_sing', 1),
    @('features_date', 'user_love', 1),
    @('bot_capabilities', 'bot_name', 1),
    @('prevention_general', 'bot_capabilities', 1),
    @('bot_personal_questions', 'bot_sexual', 1),
    @('bot_personal_questions', 'bot_residence', 1),
    @('covid_situation_infected', 'covid_situation', 1),
    @('cc_religion', 'comment_negative', 1),
    @('covid_worry', 'user_friend', 1),
    @('greeting_how_are_you', 'bot_personality', 1),
    @('prevention_medical_attention', 'covid_procedure_after_infection', 1),
    @('vocative_help', 'vocative_call', 1),
    @('bot_capabilities', 'greeting_goodbye', 1),
    @('covid_crisis_howlong', 'covid_duration', 1),
    @('user_scared', 'covid_worry', 1),
    @('bot_name', 'bot_personal_questions', 1),
    @('bot_real', 'comment_negative', 1),
    @('covid_meaning', 'covid_info', 1),
    @('vocative_call', 'comment_negative', 1),
    @('covid_origins', 'spread_general', 1),
    @('covid_food', 'spread_surfaces_food_objects', 1),
    @('prevention_medicine', 'covid_treatment', 1),
    @('covid_crisis_howlong', 'spread_general', 1),
    @('comment_positive', 'comment_smart', 1)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}
